$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.856.67'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.638.13'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9964'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.82%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5031'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2570'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07787'
$ws.Range("D11").Style = "Normal"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.647.30'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.261'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '1.862.12'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5420'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '0.0₅7881'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '25.929.04'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '197.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.376'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.925'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.886'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1139'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.830'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04865'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.254'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.181'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.532'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8885'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.602'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5514'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("D39").Value = '1.126.66'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.004'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.672'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8125'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  +7.51%  '
$ws.Range("D46").Value = '1.772.22'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05059'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("E51").Value = '  -0.76%  '
